# Update Metadata_File.xlsx to AEMO-formatted bases:
# insert a new "GVAO_PLCC" row (combined Division R & S line) ahead of the
# existing INDPRODLCC row, and widen column B to fit the longer label.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 33; this shifts the former rows 33-38 down to 34-39
# and copies the formatting (styles) of the row above, matching the
# original author's edit.
$ws.Rows(33).Insert()

# Populate the newly inserted row with the AEMO-combined Division R & S entry
$ws.Range("A33").Value = "GVAO_PLCC"
$ws.Range("B33").Value = "Gross Value Added >> Division R & S >> Other Services & Arts and Recreation Services "
$ws.Range("C33").Value = "`$'M"
$ws.Range("D33").Value = "Millions: 2020-21 prices "

# Column B needs to be a bit wider to fit the new, longer text label
$ws.Columns("B").ColumnWidth = 77.45

# Leave the selection where the author finished editing
$ws.Range("B33").Select() | Out-Null
